$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update results for Steel (column B) and minor floating point refinements
# in the Biomass (D6) and Other (D8) rows.
$ws.Range("B3").Value = 213.8864531954784
$ws.Range("D6").Value = 2676.985021085412
$ws.Range("D8").Value = 1389.575972385623
